$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the banned term for the iphone 12 64gb search row
$ws.Range("B2").Value = "mini watch 11"

# Move the active selection to B2, matching the saved workbook state
$ws.Range("B2").Select()
